$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new "NpcName" column right after NpcID (column A), pushing the
# existing NpcSex/NpcAge/NpcPersonality/NpcProplemType/NpcProblemInfo
# columns one place to the right.
$ws.Columns("B:B").Insert()

# Header
$ws.Range("B1").Value = "NpcName"

# NPC names for the six existing rows
$ws.Range("B2").Value = "Jack"
$ws.Range("B3").Value = "Tom"
$ws.Range("B4").Value = "Jerry"
$ws.Range("B5").Value = "Merry"
$ws.Range("B6").Value = "Guy"
$ws.Range("B7").Value = "Boy"

# Restore the active selection used in the saved workbook
$ws.Range("J21").Select() | Out-Null
